# Update "Valor Mora" values in the "ESTADO DE CUENTA" table.
# Row 16 corresponds to period 2107, its value moves from 48000 to 40000.
# Row 27 corresponds to period 2207, its value moves from 40000 to 48000.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("F16").Value = 40000
$ws.Range("F27").Value = 48000
